$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: "Functionality (Client App)" ---

# Row 7: mark "Show filters actived" filters work as COMPLETE
$ws1.Range("C7").Value = "COMPLETE"
$ws1.Range("C7").Style = "Buena"
$ws1.Range("C7").HorizontalAlignment = -4108
$ws1.Range("C7").VerticalAlignment = -4108

# Row 8: new function row - "Show description of one word" / DEVELOPING
$ws1.Range("A8").Value = "Show description of one word"
$ws1.Range("C8").Value = "DEVELOPING"
$ws1.Range("C8").Style = "60% - Accent6"
$ws1.Range("C8").HorizontalAlignment = -4108
$ws1.Range("C8").VerticalAlignment = -4108

# Column C width + default centered style
$ws1.Columns.Item(3).ColumnWidth = 18.109375

# Header row centered (A1:D1) while C1 additionally keeps header style (already has it)
$ws1.Range("A1:D1").HorizontalAlignment = -4108
$ws1.Range("A1:D1").VerticalAlignment = -4108

# Existing status cells centered (C2:C4)
$ws1.Range("C2:C4").HorizontalAlignment = -4108
$ws1.Range("C2:C4").VerticalAlignment = -4108

$ws1.Range("A6").HorizontalAlignment = -4108
$ws1.Range("A6").VerticalAlignment = -4108

$ws1.Range("E9").Select()

# --- Sheet2: "Bugs" ---
$ws2.Range("A1:B1").HorizontalAlignment = -4108
$ws2.Range("A1:B1").VerticalAlignment = -4108
$ws2.Range("A2").HorizontalAlignment = -4108
$ws2.Range("A2").VerticalAlignment = -4108
